# Generate Report for handback
# Adds a new handback-status row (the "3a6349ac-..." file) to the
# Overview sheet and to each per-language detail sheet (zh-cn, de-de),
# mirroring the pattern already used by the existing rows.

$wb = $excel.ActiveWorkbook

$mdName   = "3a6349ac-f09a-447a-bf50-d1d25085310e.md"
$hashPart = "3a6349ac-f09a-447a-bf50-d1d25085310e.b8d29fa28f8861860d3d42c2ba0bc1a27c5d244d"
$status   = "Handed back: in sync with en-US"
$reason   = "Include"

# ---------------------------------------------------------------
# Overview sheet: File Name | zh-cn status | de-de status
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$rowOv = 4
$wsOverview.Cells.Item($rowOv, 1).Value = $mdName
$wsOverview.Cells.Item($rowOv, 2).Value = $status
$wsOverview.Cells.Item($rowOv, 3).Value = $status

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($rowOv, 1), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/main/e2e/$mdName", "", "", $mdName) | Out-Null

# ---------------------------------------------------------------
# Per-language detail sheets
# ---------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; Lang = "zh-cn" },
    @{ Sheet = "de-de"; Lang = "de-de" }
)

foreach ($entry in $languages) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    $lang = $entry.Lang

    $xlfName = "$hashPart.$lang.xlf"
    if ($lang -eq "zh-cn") {
        $handoffDt  = "2016-01-19 06:42:44"
        $handbackDt = "2016-01-19 06:43:25"
    } else {
        $handoffDt  = "2016-01-19 06:42:54"
        $handbackDt = "2016-01-19 06:43:43"
    }

    $r = 4

    # A: Source File Name (hyperlink to the .md)
    $ws.Cells.Item($r, 1).Value = $mdName
    # B: Status
    $ws.Cells.Item($r, 2).Value = $status
    # C: Correspond Handoff File (hyperlink to the .xlf)
    $ws.Cells.Item($r, 3).Value = $xlfName
    # D: Correspond Handoff Datetime
    $ws.Cells.Item($r, 4).Value = $handoffDt
    $ws.Cells.Item($r, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    # E: Target File (hyperlink to the .md)
    $ws.Cells.Item($r, 5).Value = $mdName
    # F: Correspond Handback File (hyperlink to the .xlf)
    $ws.Cells.Item($r, 6).Value = $xlfName
    # G: Correspond Handback DateTime
    $ws.Cells.Item($r, 7).Value = $handbackDt
    # H: Handoff Reason
    $ws.Cells.Item($r, 8).Value = $reason

    $mdUrl  = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/main/e2e/$mdName"
    $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/main/ol-handback/OpenLocalizationTestOrg/oltest.$lang/xinjiang/$xlfName"

    $ws.Hyperlinks.Add($ws.Cells.Item($r, 1), $mdUrl, "", "", $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 3), $xlfUrl, "", "", $xlfName) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 5), $mdUrl, "", "", $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $xlfUrl, "", "", $xlfName) | Out-Null
}
